$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price = column D, Volume(1h) = column E)
# Force text number format so numeric/percent-looking strings are stored as literal text,
# matching the workbook's inline-string convention.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.02%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.48%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.557"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.53%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08098"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.90%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.671"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.51%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.911"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.74%"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.300"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.70%"

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.34%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9486"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.32%"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1186"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.85%"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1900"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.65%"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09686"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.26%"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04050"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.05%"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1068"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.43%"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001283"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.77%"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005947"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.54%"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.577"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.19%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3486"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.67%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.700"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.12%"

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.67%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2593"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.49%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04333"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.02%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001241"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.24%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004628"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.85%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001233"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.16%"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.01%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02673"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.94%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05486"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.78%"

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "26.28%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007719"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.52%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1396"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.97%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002113"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.12%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009446"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.70%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007092"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.16%"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.01%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003449"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.33%"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002276"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.32%"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.01%"
